$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.015.82'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").Value = '1.828.94'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''326.10'
$ws.Range("E5").Value = '  -3.02%  '
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").Value = '''0.4641'
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("D8").Value = '''0.3869'
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("D9").Value = '''0.07867'
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("D10").Value = '''0.9600'
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '''21.87'
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.921.72'
$ws.Range("E12").Value = '  +3.21%  '
$ws.Range("D13").Value = '''5.658'
$ws.Range("E13").Value = '  -3.27%  '
$ws.Range("D14").Value = '''6.890'
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").Value = '''0.06739'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '''87.07'
$ws.Range("E16").Value = '  -0.73%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = '''1.002'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").Value = '''0.000009931'
$ws.Range("E18").Value = '  -1.92%  '
$ws.Range("D19").Value = '''16.61'
$ws.Range("E19").Value = '  -2.51%  '
$ws.Range("D20").Value = '''1.001'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '28.036.12'
$ws.Range("E21").Value = '  -2.16%  '
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("E23").Value = '  -2.67%  '
$ws.Range("D24").Value = '''2.097'
$ws.Range("E24").Value = '  -1.31%  '
$ws.Range("D25").Value = '2.118.26'
$ws.Range("E25").Value = '  +2.37%  '
$ws.Range("D26").Value = '''153.70'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '''19.12'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").Value = '''5.736'
$ws.Range("E28").Value = '  -8.91%  '
$ws.Range("D29").Value = '''1.972'
$ws.Range("E29").Value = '  -2.46%  '
$ws.Range("D30").Value = '''117.29'
$ws.Range("E30").Value = '  -0.34%  '
$ws.Range("D31").Value = '''0.9372'
$ws.Range("D32").Value = '''0.09255'
$ws.Range("E32").Value = '  -2.09%  '
$ws.Range("D33").Value = '''5.290'
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("D34").Value = '''1.314'
$ws.Range("E34").Value = '  -2.89%  '
$ws.Range("D35").Value = '''3.316'
$ws.Range("E35").Value = '  -5.50%  '
$ws.Range("D36").Value = '''0.05868'
$ws.Range("E36").Value = '  -4.56%  '
$ws.Range("D37").Value = '''0.02142'
$ws.Range("E37").Value = '  -2.48%  '
$ws.Range("D38").Value = '''1.145'
$ws.Range("E38").Value = '  -1.00%  '
$ws.Range("D39").Value = '''7.770'
$ws.Range("E39").Value = '  +2.38%  '
$ws.Range("D40").Value = '''0.5585'
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("D41").Value = '''9.875'
$ws.Range("E41").Value = '  -2.25%  '
$ws.Range("D42").Value = '''0.1758'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("D43").Value = '''1.215'
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("D44").Value = '''11.66'
$ws.Range("E44").Value = '  -1.16%  '
$ws.Range("D45").Value = '''0.5261'
$ws.Range("E45").Value = '  -2.32%  '
$ws.Range("D46").Value = '''0.07016'
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("E47").Value = '  -8.48%  '
$ws.Range("D48").Value = '''1.829'
$ws.Range("E48").Value = '  -4.30%  '
$ws.Range("D49").Value = '''112.97'
$ws.Range("E49").Value = '  -1.79%  '
$ws.Range("D50").Value = '''0.9995'
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("E51").Value = '  -0.03%  '
